$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 404.08334
$ws.Range("I8").Value = 205.44444
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 616.33332
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = -477.33332
$ws.Range("N8").Value = -3278

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2145.8547
$ws.Range("I138").Value = 2171.15
$ws.Range("J138").Value = 2133.8096
$ws.Range("K138").Value = 6513.450000000001
$ws.Range("L138").Value = 6401.4288
$ws.Range("M138").Value = -1373.450000000001
$ws.Range("N138").Value = -16681.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4442.222
$ws.Range("I45").Value = 3995
$ws.Range("J45").Value = 4800
$ws.Range("K45").Value = 3995
$ws.Range("L45").Value = 4800
$ws.Range("M45").Value = -3618
$ws.Range("N45").Value = -5554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101622

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -308112

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 33000
$ws.Range("J82").Value = 33000
$ws.Range("L82").Value = 33000
$ws.Range("N82").Value = -33722

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H85").Value = 33000
$ws.Range("J85").Value = 33000
$ws.Range("L85").Value = 33000
$ws.Range("N85").Value = -35496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 32287096
$ws.Range("J86").Value = 32287096
$ws.Range("L86").Value = 32287096
$ws.Range("N86").Value = -32289468

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 32287096
$ws.Range("J89").Value = 32287096
$ws.Range("L89").Value = 96861288
$ws.Range("N89").Value = -96873144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -117480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 73929
$ws.Range("J141").Value = 73929
$ws.Range("L141").Value = 73929
$ws.Range("N141").Value = -84289

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101622

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 46250
$ws.Range("J69").Value = 46250
$ws.Range("L69").Value = 46250
$ws.Range("N69").Value = -47872

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -308112

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H72").Value = 46250
$ws.Range("J72").Value = 46250
$ws.Range("L72").Value = 138750
$ws.Range("N72").Value = -146862

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 678.6111
$ws.Range("I94").Value = 636.0714
$ws.Range("J94").Value = 827.5
$ws.Range("K94").Value = 636.0714
$ws.Range("L94").Value = 827.5
$ws.Range("M94").Value = -185.0714
$ws.Range("N94").Value = -1729.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2204
$ws.Range("I99").Value = 2307.5
$ws.Range("K99").Value = 2307.5
$ws.Range("M99").Value = -809.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 64937.8
$ws.Range("J100").Value = 64937.8
$ws.Range("L100").Value = 64937.8
$ws.Range("N100").Value = -67101.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 30000
$ws.Range("J88").Value = 30000
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30812

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 30000
$ws.Range("J91").Value = 30000
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -32808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1895.4286
$ws.Range("I134").Value = 1501.6
$ws.Range("J134").Value = 2880
$ws.Range("K134").Value = 4504.799999999999
$ws.Range("L134").Value = 8640
$ws.Range("M134").Value = -1969.799999999999
$ws.Range("N134").Value = -13710

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 979.92
$ws.Range("J113").Value = 1089.9
$ws.Range("L113").Value = 3269.7
$ws.Range("N113").Value = -7609.700000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2575.125
$ws.Range("I97").Value = 2096.5
$ws.Range("K97").Value = 2096.5
$ws.Range("M97").Value = -1600.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 32999.5
$ws.Range("I74").Value = 32999
$ws.Range("J74").Value = 33000
$ws.Range("K74").Value = 32999
$ws.Range("L74").Value = 33000
$ws.Range("M74").Value = -32001
$ws.Range("N74").Value = -34996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 32999.5
$ws.Range("I77").Value = 32999
$ws.Range("J77").Value = 33000
$ws.Range("K77").Value = 98997
$ws.Range("L77").Value = 99000
$ws.Range("M77").Value = -94005
$ws.Range("N77").Value = -108984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 59999.5
$ws.Range("J140").Value = 59999.5
$ws.Range("L140").Value = 59999.5
$ws.Range("N140").Value = -70359.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 1000000000
$ws.Range("J32").Value = 1000000000
$ws.Range("L32").Value = 1000000000
$ws.Range("N32").Value = -1000000634

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13935
$ws.Range("I81").Value = 13935
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 27870
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -26809
$ws.Range("N81").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 13935
$ws.Range("I84").Value = 13935
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 139350
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -134046
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 60330
$ws.Range("J94").Value = 60330
$ws.Range("L94").Value = 60330
$ws.Range("N94").Value = -62132

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 961.2
$ws.Range("I107").Value = 935.3333
$ws.Range("K107").Value = 2805.9999
$ws.Range("M107").Value = -885.9998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 738.6
$ws.Range("I113").Value = 818.86365
$ws.Range("J113").Value = 150
$ws.Range("K113").Value = 2456.59095
$ws.Range("L113").Value = 450
$ws.Range("M113").Value = -286.5909499999998
$ws.Range("N113").Value = -4790

